# "Örnek 31 - Eğer Örneği" — add the IF() height-classification formulas
# and fill in the "Numara / Ad Soyad / Bölüm" info box.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Threshold value used by the IF() formulas below.
$ws.Range("D2").Value = 169

# D4 gets its own (non-shared) formula...
$ws.Range("D4").Formula = '=IF(C4>$D$2,"UZUN","KISA")'

# ...and D5:D11 is filled with the same formula (creates the D5:D11 shared
# formula group seen in the saved file).
$ws.Range("D5:D11").Formula = '=IF(C5>$D$2,"UZUN","KISA")'

# Student info box (H6:H8 labels already exist: Numara:/Ad Soyad:/Bölüm:).
$ws.Range("I6").Value = 20215070019
$ws.Range("I7").Value = "KÜBRA ÇABUK"
$ws.Range("I8").Value = "YBS"

# Match the saved selection state.
$null = $ws.Range("J11").Select()
